$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column -> new text value, derived from the commit diff.
# Columns D, E, G hold numeric-looking / percent-looking text and must be
# written with a leading apostrophe so Excel keeps them as text (matching
# the original inlineStr cells) instead of coercing them to numbers.
# Column F holds a "D-M-YYYY" style date string that Excel does not
# auto-recognise as a date, so it is safe to assign directly.
$updates = @{
    2 = @{ "D" = "293.61"; "E" = "1.75%"; "F" = "20-1-2023"; "G" = "0" }
    3 = @{ "D" = "31.04"; "E" = "1.31%"; "F" = "20-1-2023"; "G" = "0" }
    4 = @{ "D" = "4.905"; "E" = "-0.32%"; "F" = "20-1-2023"; "G" = "0" }
    5 = @{ "D" = "0.07371"; "E" = "3.71%"; "F" = "20-1-2023"; "G" = "0" }
    6 = @{ "D" = "2.274"; "E" = "27.03%"; "F" = "20-1-2023"; "G" = "0" }
    7 = @{ "D" = "7.688"; "E" = "0.03%"; "F" = "20-1-2023"; "G" = "0" }
    8 = @{ "D" = "3.751"; "F" = "20-1-2023"; "G" = "0" }
    9 = @{ "D" = "0.9152"; "E" = "2.02%"; "F" = "20-1-2023"; "G" = "0" }
    10 = @{ "D" = "0.1693"; "E" = "3.17%"; "F" = "20-1-2023"; "G" = "0" }
    11 = @{ "D" = "0.08424"; "E" = "11.26%"; "F" = "20-1-2023"; "G" = "0" }
    12 = @{ "D" = "0.08317"; "E" = "3.22%"; "F" = "20-1-2023"; "G" = "0" }
    13 = @{ "D" = "0.03078"; "E" = "0.63%"; "F" = "20-1-2023"; "G" = "0" }
    14 = @{ "D" = "0.1007"; "E" = "0.72%"; "F" = "20-1-2023"; "G" = "0" }
    15 = @{ "D" = "0.001519"; "E" = "0.90%"; "F" = "20-1-2023"; "G" = "0" }
    16 = @{ "D" = "0.005754"; "E" = "0.44%"; "F" = "20-1-2023"; "G" = "0" }
    17 = @{ "D" = "3.482"; "E" = "0.51%"; "F" = "20-1-2023"; "G" = "0" }
    18 = @{ "D" = "2.077"; "E" = "-1.02%"; "F" = "20-1-2023"; "G" = "0" }
    19 = @{ "E" = "1.68%"; "F" = "20-1-2023"; "G" = "0" }
    20 = @{ "F" = "20-1-2023"; "G" = "0" }
    21 = @{ "D" = "3.980"; "E" = "-1.58%"; "F" = "20-1-2023"; "G" = "0" }
    22 = @{ "D" = "0.2099"; "E" = "4.88%"; "F" = "20-1-2023"; "G" = "0" }
    23 = @{ "D" = "0.04551"; "E" = "1.25%"; "F" = "20-1-2023"; "G" = "0" }
    24 = @{ "D" = "0.001210"; "E" = "-0.22%"; "F" = "20-1-2023"; "G" = "0" }
    25 = @{ "D" = "0.004332"; "E" = "-6.63%"; "F" = "20-1-2023"; "G" = "0" }
    26 = @{ "D" = "0.0001302"; "E" = "3.98%"; "F" = "20-1-2023"; "G" = "0" }
    27 = @{ "D" = "0.0003393"; "F" = "20-1-2023"; "G" = "0" }
    28 = @{ "F" = "20-1-2023"; "G" = "0" }
    29 = @{ "F" = "20-1-2023"; "G" = "0" }
    30 = @{ "F" = "20-1-2023"; "G" = "0" }
    31 = @{ "F" = "20-1-2023"; "G" = "0" }
    32 = @{ "F" = "20-1-2023"; "G" = "0" }
    33 = @{ "F" = "20-1-2023"; "G" = "0" }
    34 = @{ "F" = "20-1-2023"; "G" = "0" }
    35 = @{ "F" = "20-1-2023"; "G" = "0" }
    36 = @{ "F" = "20-1-2023"; "G" = "0" }
    37 = @{ "F" = "20-1-2023"; "G" = "0" }
    38 = @{ "F" = "20-1-2023"; "G" = "0" }
    39 = @{ "D" = "0.01610"; "E" = "-0.24%"; "F" = "20-1-2023"; "G" = "0" }
    40 = @{ "D" = "0.04452"; "E" = "2.34%"; "F" = "20-1-2023"; "G" = "0" }
    41 = @{ "D" = "0.007348"; "E" = "0.03%"; "F" = "20-1-2023"; "G" = "0" }
    42 = @{ "D" = "0.008943"; "F" = "20-1-2023"; "G" = "0" }
    43 = @{ "D" = "0.1327"; "E" = "1.94%"; "F" = "20-1-2023"; "G" = "0" }
    44 = @{ "D" = "0.001979"; "E" = "-1.23%"; "F" = "20-1-2023"; "G" = "0" }
    45 = @{ "D" = "0.009172"; "E" = "-5.53%"; "F" = "20-1-2023"; "G" = "0" }
    46 = @{ "D" = "0.00006027"; "E" = "2.50%"; "F" = "20-1-2023"; "G" = "0" }
    47 = @{ "D" = "0.00000000751"; "E" = "-0.02%"; "F" = "20-1-2023"; "G" = "0" }
    48 = @{ "E" = "0.47%"; "F" = "20-1-2023"; "G" = "0" }
    49 = @{ "E" = "-3.50%"; "F" = "20-1-2023"; "G" = "0" }
    50 = @{ "D" = "0.00002103"; "E" = "-0.02%"; "F" = "20-1-2023"; "G" = "0" }
    51 = @{ "D" = "0.0002003"; "E" = "-0.02%"; "F" = "20-1-2023"; "G" = "0" }
}

foreach ($row in $updates.Keys) {
    $rowData = $updates[$row]
    foreach ($col in $rowData.Keys) {
        $value = $rowData[$col]
        $cellRef = "$col$row"
        if ($col -eq "F") {
            $ws.Range($cellRef).Value = $value
        } else {
            $ws.Range($cellRef).Value = "'" + $value
        }
    }
}
